$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (reflected in workbook.xml <sheet name=.../>)
$ws.Name = "Comprar test excel fail"

# "Samsung Galaxy Tab 10.1" (row 2) is replaced by "fail", and its cart
# status changes from "In Cart" to "Not In Cart"
$ws.Range("A2").Value = "fail"
$ws.Range("B2").Value = "Not In Cart"

# "Nikon D300" (row 4) is also replaced by "fail" / "Not In Cart"
$ws.Range("A4").Value = "fail"
$ws.Range("B4").Value = "Not In Cart"
